# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary header updates ---
# VALOR MORA total (E11)
$ws.Range("E11").Value = 1621640
# Cant. Trabajadores (C13) : 3 -> 2 (HEBERT removed)
$ws.Range("C13").Value = 2
# Cant. Periodos (F13) : 31 -> 32
$ws.Range("F13").Value = 32

# --- Data table (rows 16-53) ---
# Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora, F=Valor Mora, G=Salario Basico
# Worker: AMAURY CABEZA COHEN (CC 1143333216) now covers periods 2301..2508
# Worker: JAIRO ALONSO JUNIOR HERRERA LEON (CC 1047459087) covers periods 2502,2503,2505,2506,2507,2508
# Worker HEBERT EDUARDO WATTS MATURANA removed entirely.

$rows = @(
    @(16, "CC", "1143333216", "AMAURY CABEZA COHEN", "2301", 40000, 1000000),
    @(17, "CC", "1143333216", "AMAURY CABEZA COHEN", "2302", 40000, 1000000),
    @(18, "CC", "1143333216", "AMAURY CABEZA COHEN", "2303", 40000, 1000000),
    @(19, "CC", "1143333216", "AMAURY CABEZA COHEN", "2304", 40000, 1000000),
    @(20, "CC", "1143333216", "AMAURY CABEZA COHEN", "2305", 40000, 1000000),
    @(21, "CC", "1143333216", "AMAURY CABEZA COHEN", "2306", 40000, 1000000),
    @(22, "CC", "1143333216", "AMAURY CABEZA COHEN", "2307", 40000, 1000000),
    @(23, "CC", "1143333216", "AMAURY CABEZA COHEN", "2308", 40000, 1000000),
    @(24, "CC", "1143333216", "AMAURY CABEZA COHEN", "2309", 40000, 1000000),
    @(25, "CC", "1143333216", "AMAURY CABEZA COHEN", "2310", 40000, 1000000),
    @(26, "CC", "1143333216", "AMAURY CABEZA COHEN", "2311", 40000, 1000000),
    @(27, "CC", "1143333216", "AMAURY CABEZA COHEN", "2312", 40000, 1000000),
    @(28, "CC", "1143333216", "AMAURY CABEZA COHEN", "2401", 40000, 1000000),
    @(29, "CC", "1143333216", "AMAURY CABEZA COHEN", "2402", 40000, 1000000),
    @(30, "CC", "1143333216", "AMAURY CABEZA COHEN", "2403", 40000, 1000000),
    @(31, "CC", "1143333216", "AMAURY CABEZA COHEN", "2404", 40000, 1000000),
    @(32, "CC", "1143333216", "AMAURY CABEZA COHEN", "2405", 40000, 1000000),
    @(33, "CC", "1143333216", "AMAURY CABEZA COHEN", "2406", 40000, 1000000),
    @(34, "CC", "1143333216", "AMAURY CABEZA COHEN", "2407", 40000, 1000000),
    @(35, "CC", "1143333216", "AMAURY CABEZA COHEN", "2408", 40000, 1000000),
    @(36, "CC", "1143333216", "AMAURY CABEZA COHEN", "2409", 40000, 1000000),
    @(37, "CC", "1143333216", "AMAURY CABEZA COHEN", "2410", 40000, 1000000),
    @(38, "CC", "1143333216", "AMAURY CABEZA COHEN", "2411", 40000, 1000000),
    @(39, "CC", "1143333216", "AMAURY CABEZA COHEN", "2412", 40000, 1000000),
    @(40, "CC", "1143333216", "AMAURY CABEZA COHEN", "2501", 40000, 1000000),
    @(41, "CC", "1143333216", "AMAURY CABEZA COHEN", "2502", 40000, 1000000),
    @(42, "CC", "1047459087", "JAIRO ALONSO JUNIOR HERRERA LEON", "2502", 56940, 1423500),
    @(43, "CC", "1143333216", "AMAURY CABEZA COHEN", "2503", 40000, 1000000),
    @(44, "CC", "1047459087", "JAIRO ALONSO JUNIOR HERRERA LEON", "2503", 56940, 1423500),
    @(45, "CC", "1143333216", "AMAURY CABEZA COHEN", "2504", 40000, 1000000),
    @(46, "CC", "1143333216", "AMAURY CABEZA COHEN", "2505", 40000, 1000000),
    @(47, "CC", "1047459087", "JAIRO ALONSO JUNIOR HERRERA LEON", "2505", 56940, 1423500),
    @(48, "CC", "1143333216", "AMAURY CABEZA COHEN", "2506", 40000, 1000000),
    @(49, "CC", "1047459087", "JAIRO ALONSO JUNIOR HERRERA LEON", "2506", 56940, 1423500),
    @(50, "CC", "1143333216", "AMAURY CABEZA COHEN", "2507", 40000, 1000000),
    @(51, "CC", "1047459087", "JAIRO ALONSO JUNIOR HERRERA LEON", "2507", 56940, 1423500),
    @(52, "CC", "1143333216", "AMAURY CABEZA COHEN", "2508", 40000, 1000000),
    @(53, "CC", "1047459087", "JAIRO ALONSO JUNIOR HERRERA LEON", "2508", 56940, 1423500)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}

Write-Host "Edit complete"
